$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.646.54'
$ws.Range("E2").Value = '  +0.58%  '

$ws.Range("D3").Value = '1.830.54'
$ws.Range("E3").Value = '  +1.19%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.006'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4683'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.42%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3601'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07150'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.75%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9320'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.70%  '

$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.47'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.31%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07650'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.31%  '

$ws.Range("D13").Value = '1.862.82'
$ws.Range("E13").Value = '  +3.04%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.262'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.71%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.347'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.15%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.63'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.81%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.008'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.21%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008549'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.42%  '

$ws.Range("E19").Value = '  +0.22%  '

$ws.Range("D20").Value = '26.665.05'
$ws.Range("E20").Value = '  +0.50%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.019'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.78%  '

$ws.Range("D23").Value = '2.072.56'
$ws.Range("E23").Value = '  +1.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.57'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.910'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.78'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.56%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.94'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.56%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.996'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.99%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.59'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.873'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.26%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08827'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.33%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.165'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.38%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.856'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.67%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.166'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.81%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7380'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.443'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.076'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.17%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01924'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.75%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.949'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.04%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05156'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.96%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.908'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.76%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5068'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.07%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1498'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.07%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.117'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.07%  '

$ws.Range("E45").Value = '  +0.35%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4653'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.45%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.08'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.21%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '98.70'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.576'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.07%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06026'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.74%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.86'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.14%  '
